# FIX VALVE RATING FOR VALVE COMPARE
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (item 1): BAV21A0F1C -> BAV24G0I1C ---
$ws.Range("B2").Value = "BAV24G0I1C"
$ws.Range("C2").Value = "BALL VALVE W/INTEGRAL WELDED 2 NIPPLES, FB, FLOATING BALL, API 608, API 598, A105, CL 800, SW W/2 PE NIPPLES, MNF STD, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, API 607, LO"
$ws.Range("D2").Value = "0,5"
$ws.Range("E2").Value = "1,00"

# --- Row 3 (item 2): BAV21A0F1C -> BAV24G0I1C ---
$ws.Range("B3").Value = "BAV24G0I1C"
$ws.Range("C3").Value = "BALL VALVE W/INTEGRAL WELDED 2 NIPPLES, FB, FLOATING BALL, API 608, API 598, A105, CL 800, SW W/2 PE NIPPLES, MNF STD, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, API 607, LO"
$ws.Range("D3").Value = "0,75"
$ws.Range("E3").Value = "1,00"
$ws.Range("G3").Value = "CSO"

# --- Row 4 (item 3): BAV21A0I1C -> BAV24G0I1C ---
$ws.Range("B4").Value = "BAV24G0I1C"
$ws.Range("C4").Value = "BALL VALVE W/INTEGRAL WELDED 2 NIPPLES, FB, FLOATING BALL, API 608, API 598, A105, CL 800, SW W/2 PE NIPPLES, MNF STD, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, API 607, LO"
$ws.Range("D4").Value = "1,00"

# --- Row 5 (item 4): BAV21A0I1C -> BAV24G0I1C ---
$ws.Range("B5").Value = "BAV24G0I1C"
$ws.Range("C5").Value = "BALL VALVE W/INTEGRAL WELDED 2 NIPPLES, FB, FLOATING BALL, API 608, API 598, A105, CL 800, SW W/2 PE NIPPLES, MNF STD, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, API 607, LO"
$ws.Range("D5").Value = "1,00"
$ws.Range("E5").Value = "1,00"
$ws.Range("G5").Value = "CSO"

# --- Row 6 (item 5): BAV24G0I1C -> CKV21A0B2B ---
$ws.Range("B6").Value = "CKV21A0B2B"
$ws.Range("C6").Value = "SWING CHECK VALVE FL, API 594, API 598, A216 GR.WCB, CL 150, INST HORIZ/VERT, RF, B16.5, BOLTED COVER, SPW SS304/GRAPH, RENEWABLE SEATS, TRIM #8"
$ws.Range("D6").Value = "3,00"

# --- Row 7 (item 6): BAV24G0I1C -> CLV24F0B2B ---
$ws.Range("B7").Value = "CLV24F0B2B"
$ws.Range("C7").Value = "LIFT CHECK VALVE SW, API 602, A105, CL 800, INST HORIZ/VERT, SW, B16.11, BOLTED COVER, SPW SS304/GRAPH, PISTON TYPE OBTURATOR, TRIM #8"
$ws.Range("D7").Value = "0,5"
$ws.Range("E7").Value = "1,00"

# --- Row 8 (item 7): BAV24G0I1C -> GAV24F0B2B ---
$ws.Range("B8").Value = "GAV24F0B2B"
$ws.Range("C8").Value = "GATE VALVE SW, API 602, API 598, A105, CL 800, SW, B16.11, BB, SPW SS304/GRAPH, PKG GRAPH, TRIM #8, RENEWABLE SEATS, SOLID WEDGE, STEM OS&Y/RSNRO, HO"
$ws.Range("D8").Value = "0,75"
$ws.Range("E8").Value = "2,00"

# --- Row 9 (item 8): CKV21A0B2B -> MFV21A0I2I ---
$ws.Range("B9").Value = "MFV21A0I2I"
$ws.Range("C9").Value = "INTEGRAL MONO FLANGE DBB NEEDLE MULTI-VALVE, EEMUA 182, A105, CL 150, RF/NPTF, B16.5 AND B1.20.1, BB, SPW SS304/GRAPH, PKG GRAPH; SS316 STEM, SEATS&STEM TIP, S, SWIVEL NEEDLE, STEM OS&Y/RSRO, T-HANDLE"
$ws.Range("D9").Value = "0,75"
$ws.Range("E9").Value = "6,00"

# --- Remove former rows 10-16 (items 9-15), no longer present after the fix ---
$ws.Rows("10:16").Delete()
